$wb = $excel.ActiveWorkbook

# Map of sheet name -> cell -> new value for the "Elapsed Duration(Hrs)" column (G)
$updates = @{
    "R1" = @{ "G2" = "3926:24:33"; "G3" = "65:57:11" }
    "R2" = @{ "G2" = "12107:48:11"; "G3" = "3237:31:40"; "G4" = "475:43:14" }
    "R4" = @{ "G2" = "2953:38:00"; "G3" = "180:50:15" }
    "R5" = @{ "G2" = "427:36:59" }
    "R6" = @{ "G2" = "68:09:17" }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellAddr in $cellUpdates.Keys) {
        $ws.Range($cellAddr).Value = $cellUpdates[$cellAddr]
    }
}
